$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 previously held a free-text answer in C11 ("I have USD, PHL, THB and
# VND accounts"). The survey got normalized to two Yes/No columns
# (CurrenciesOutUS / wantUSCurrency), so C11 becomes "Yes" and a new D11
# "Yes" is added alongside it.
$ws.Range("C11").Value = "Yes"
$ws.Range("C11").Font.Name = "Calibri"

$ws.Range("D11").Value = "Yes"
$ws.Range("D11").Font.Name = "Calibri"
$ws.Range("D11").IndentLevel = 0

# Leave the cursor where the edit was made.
$ws.Range("D19").Select()
